# Cambios de los excel e icono de usuarios
#
# The "Carrera" (major/career) value "Mecatronica" was renamed to "Sistemas"
# for the two data rows, and the last-selected cell on the sheet changed
# to E7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Sistemas"
$ws.Range("D3").Value = "Sistemas"

$ws.Range("E7").Select()
